$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# "Content Placeholder 2" holds the bullet text we're appending a new
# paragraph to (the title shape is Shapes.Item(1)).
$shape = $s.Shapes.Item("Content Placeholder 2")
$tr = $shape.TextFrame.TextRange

# Typing a new line at the end of the existing text, then the first
# sentence of the new paragraph.
$run1 = $tr.InsertAfter([char]13 + "The data for this project have been sourced from ")

# "kagel" is typed as a separate word/run right after.
$full = $shape.TextFrame.TextRange
$run2 = $full.InsertAfter("kagel")
